$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D stay text (match source formatting),
# matching the original inline-string cells which are never coerced to real numbers.

$ws.Range("D2").Value = "30.141.37"
$ws.Range("E2").Value = "  -1.80%  "
$ws.Range("D3").Value = "1.830.72"
$ws.Range("E3").Value = "  -3.40%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.64"
$ws.Range("E5").Value = "  -3.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4642"
$ws.Range("E7").Value = "  -4.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2694"
$ws.Range("E8").Value = "  -6.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06256"
$ws.Range("E9").Value = "  -4.68%  "
$ws.Range("D10").Value = "1.817.78"
$ws.Range("E10").Value = "  -3.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07363"
$ws.Range("E11").Value = "  -1.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.03"
$ws.Range("E12").Value = "  -5.13%  "
$ws.Range("E13").Value = "  -4.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "83.06"
$ws.Range("E14").Value = "  -5.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6188"
$ws.Range("E15").Value = "  -7.66%  "
$ws.Range("D16").Value = "30.073.06"
$ws.Range("E16").Value = "  -1.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9996"
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.19"
$ws.Range("E18").Value = "  -1.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007264"
$ws.Range("E19").Value = "  -4.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.39"
$ws.Range("E20").Value = "  -6.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9992"
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("D22").Value = "2.067.90"
$ws.Range("E22").Value = "  -3.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.829"
$ws.Range("E23").Value = "  -8.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.838"
$ws.Range("E24").Value = "  -6.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "164.73"
$ws.Range("E25").Value = "  -3.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.089"
$ws.Range("E26").Value = "  -3.15%  "
$ws.Range("E27").Value = "  -6.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.838"
$ws.Range("E28").Value = "  -6.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1016"
$ws.Range("E29").Value = "  -1.26%  "
$ws.Range("E30").Value = "  -2.03%  "
$ws.Range("E31").Value = "  -7.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.751"
$ws.Range("E32").Value = "  -6.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04785"
$ws.Range("E33").Value = "  -5.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.124"
$ws.Range("E34").Value = "  -7.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6988"
$ws.Range("E35").Value = "  -7.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.687"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01809"
$ws.Range("E37").Value = "  -3.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.606"
$ws.Range("E38").Value = "  -1.71%  "
$ws.Range("E39").Value = "  -3.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.918"
$ws.Range("E40").Value = "  -7.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9999"
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "102.31"
$ws.Range("E42").Value = "  -4.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.492"
$ws.Range("E43").Value = "  -3.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3990"
$ws.Range("E44").Value = "  -7.41%  "
$ws.Range("E45").Value = "  -7.08%  "
$ws.Range("E46").Value = "  -7.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "59.39"
$ws.Range("E47").Value = "  -7.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.392"
$ws.Range("E48").Value = "  -6.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05523"
$ws.Range("E49").Value = "  -2.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "32.46"
$ws.Range("E50").Value = "  -4.76%  "
$ws.Range("E51").Value = "  -9.61%  "
